# Helper: set a cell's value as literal text (avoids Excel's automatic
# number/date coercion of numeric-looking strings such as "10.00" or "012159").
function Set-TextValue {
    param($Range, [string]$Text)
    $escaped = $Text -replace '"', '""'
    $Range.Formula = '="' + $escaped + '"'
    $Range.Copy($null)
    $Range.PasteSpecial(-4163)  # xlPasteValues
}

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)        # "总计"
$q2sheet = $wb.Worksheets.Item(2)      # "2022-Q2" (will become 2022-Q3's template)

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    right under the header, pushing the existing rows down.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# Restore the bordered/bold/centered style on A2 (same as the other
# rows in column A) by copying formats from the row that used to be
# row 2 and is now row 3.
$total.Range("A3").Copy($null)
$total.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# The Insert() operation also carried over formatting into B2:D2 from
# the header row above; clear that so they match the plain data rows.
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 1.9

# ---------------------------------------------------------------------
# 2. Add the new "2022-Q3" sheet (holdings detail), positioned right
#    after "总计". We clone the existing "2022-Q2" sheet so the new
#    sheet inherits identical layout/formatting/column types, then
#    overwrite the values.
# ---------------------------------------------------------------------
$q2sheet.Copy($null, $total)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Row 2
Set-TextValue $newSheet.Range("B2") "012159"
$newSheet.Range("C2").Value = "财通资管健康产业混合A"
Set-TextValue $newSheet.Range("D2") "10.00"
Set-TextValue $newSheet.Range("E2") "94.52"
Set-TextValue $newSheet.Range("F2") "9.91"
Set-TextValue $newSheet.Range("G2") "0.9910"
$newSheet.Range("H2").Value = 1

# Row 3
Set-TextValue $newSheet.Range("B3") "519087"
$newSheet.Range("C3").Value = "新华优选分红混合"
Set-TextValue $newSheet.Range("D3") "10.82"
Set-TextValue $newSheet.Range("E3") "89.55"
Set-TextValue $newSheet.Range("F3") "3.39"
Set-TextValue $newSheet.Range("G3") "0.3668"
$newSheet.Range("H3").Value = 10

# Row 4
Set-TextValue $newSheet.Range("B4") "001040"
$newSheet.Range("C4").Value = "新华策略精选股票"
Set-TextValue $newSheet.Range("D4") "6.78"
Set-TextValue $newSheet.Range("E4") "94.54"
Set-TextValue $newSheet.Range("F4") "4.04"
Set-TextValue $newSheet.Range("G4") "0.2739"
$newSheet.Range("H4").Value = 10

# Row 5
Set-TextValue $newSheet.Range("B5") "012160"
$newSheet.Range("C5").Value = "财通资管健康产业混合C"
Set-TextValue $newSheet.Range("D5") "1.95"
Set-TextValue $newSheet.Range("E5") "94.52"
Set-TextValue $newSheet.Range("F5") "9.91"
Set-TextValue $newSheet.Range("G5") "0.1932"
$newSheet.Range("H5").Value = 1

# Row 6
Set-TextValue $newSheet.Range("B6") "005044"
$newSheet.Range("C6").Value = "国寿安保健康科学混合C"
Set-TextValue $newSheet.Range("D6") "0.75"
Set-TextValue $newSheet.Range("E6") "87.70"
Set-TextValue $newSheet.Range("F6") "4.65"
Set-TextValue $newSheet.Range("G6") "0.0349"
$newSheet.Range("H6").Value = 4

# Row 7
Set-TextValue $newSheet.Range("B7") "001294"
$newSheet.Range("C7").Value = "新华战略新兴产业灵活配置混合"
Set-TextValue $newSheet.Range("D7") "0.99"
Set-TextValue $newSheet.Range("E7") "93.49"
Set-TextValue $newSheet.Range("F7") "3.44"
Set-TextValue $newSheet.Range("G7") "0.0341"
$newSheet.Range("H7").Value = 8

# Row 8
Set-TextValue $newSheet.Range("B8") "005043"
$newSheet.Range("C8").Value = "国寿安保健康科学混合A"
Set-TextValue $newSheet.Range("D8") "0.08"
Set-TextValue $newSheet.Range("E8") "87.70"
Set-TextValue $newSheet.Range("F8") "4.65"
Set-TextValue $newSheet.Range("G8") "0.0037"
$newSheet.Range("H8").Value = 4

# The template sheet (old "2022-Q2") had 12 data rows (through row 12);
# the new "2022-Q3" sheet only needs 7 (through row 8). Remove the
# leftover rows 9-12 copied from the template.
$newSheet.Range("A9:H12").ClearContents()
$newSheet.Range("A9:H12").ClearFormats()
